# 9.2.1 sheet: add a 2022 data column (S) and revise the 2019-2021 figures
# for both indicator rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S, inheriting formatting (styles, borders, etc.)
# from the column immediately to its left (R) - this gives the new column
# the same per-row styles already used for the existing year columns.
$ws.Columns("S").Insert()

# New 2022 year header
$ws.Range("S3").Value = 2022

# Row 4 (GVA share of manufacturing output in GDP, %): revise 2019-2021,
# add 2022
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8
$ws.Range("S4").Value = 13.6

# Row 5 (GVA of manufacturing industry in GDP per capita): revise 2019-2021,
# add 2022
$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5
$ws.Range("S5").Value = 20

# Match the author's final selection on the new column's header-row cell
[void]$ws.Range("S2").Select()
